$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,4).Value = '24.296.62'
$ws.Cells.Item(2,5).Value = '  +14.08%  '

# Row 3
$ws.Cells.Item(3,4).Value = '1.680.23'
$ws.Cells.Item(3,5).Value = '  +8.74%  '

# Row 4
$ws.Cells.Item(4,5).Value = '  -0.35%  '

# Row 5
$ws.Cells.Item(5,4).Value = '''307.43'
$ws.Cells.Item(5,5).Value = '  +8.74%  '

# Row 6
$ws.Cells.Item(6,4).Value = '''0.9961'
$ws.Cells.Item(6,5).Value = '  +3.84%  '

# Row 7
$ws.Cells.Item(7,4).Value = '''0.3725'
$ws.Cells.Item(7,5).Value = '  +2.44%  '

# Row 8
$ws.Cells.Item(8,4).Value = '''0.3443'
$ws.Cells.Item(8,5).Value = '  +7.75%  '

# Row 9
$ws.Cells.Item(9,5).Value = '  +17.59%  '

# Row 10
$ws.Cells.Item(10,4).Value = '''1.190'
$ws.Cells.Item(10,5).Value = '  +8.37%  '

# Row 11
$ws.Cells.Item(11,4).Value = '''0.07300'
$ws.Cells.Item(11,5).Value = '  +6.91%  '

# Row 12
$ws.Cells.Item(12,4).Value = '''0.9978'
$ws.Cells.Item(12,5).Value = '  -0.04%  '

# Row 13
$ws.Cells.Item(13,4).Value = '''20.66'
$ws.Cells.Item(13,5).Value = '  +9.76%  '

# Row 14
$ws.Cells.Item(14,4).Value = '''6.116'
$ws.Cells.Item(14,5).Value = '  +7.52%  '

# Row 15
$ws.Cells.Item(15,4).Value = '''6.783'
$ws.Cells.Item(15,5).Value = '  +6.47%  '

# Row 16
$ws.Cells.Item(16,4).Value = '1.676.39'
$ws.Cells.Item(16,5).Value = '  +9.39%  '

# Row 17
$ws.Cells.Item(17,4).Value = '''0.00001113'
$ws.Cells.Item(17,5).Value = '  +5.78%  '

# Row 18
$ws.Cells.Item(18,4).Value = '''0.9963'
$ws.Cells.Item(18,5).Value = '  +3.86%  '

# Row 19
$ws.Cells.Item(19,4).Value = '''0.06712'
$ws.Cells.Item(19,5).Value = '  +10.54%  '

# Row 20
$ws.Cells.Item(20,4).Value = '''81.84'
$ws.Cells.Item(20,5).Value = '  +12.68%  '

# Row 21
$ws.Cells.Item(21,4).Value = '''16.51'
$ws.Cells.Item(21,5).Value = '  +9.57%  '

# Row 22
$ws.Cells.Item(22,4).Value = '''6.146'
$ws.Cells.Item(22,5).Value = '  +7.73%  '

# Row 23
$ws.Cells.Item(23,4).Value = '''12.03'
$ws.Cells.Item(23,5).Value = '  +5.75%  '

# Row 24
$ws.Cells.Item(24,4).Value = '24.232.55'
$ws.Cells.Item(24,5).Value = '  +13.49%  '

# Row 25
$ws.Cells.Item(25,4).Value = '''2.403'
$ws.Cells.Item(25,5).Value = '  +3.93%  '

# Row 26
$ws.Cells.Item(26,4).Value = '''2.676'
$ws.Cells.Item(26,5).Value = '  +20.39%  '

# Row 27
$ws.Cells.Item(27,4).Value = '''3.360'
$ws.Cells.Item(27,5).Value = '  -9.58%  '

# Row 28
$ws.Cells.Item(28,4).Value = '''151.74'
$ws.Cells.Item(28,5).Value = '  +2.31%  '

# Row 29
$ws.Cells.Item(29,4).Value = '''19.54'
$ws.Cells.Item(29,5).Value = '  +10.31%  '

# Row 30
$ws.Cells.Item(30,4).Value = '1.858.74'
$ws.Cells.Item(30,5).Value = '  +9.13%  '

# Row 31
$ws.Cells.Item(31,4).Value = '''127.31'
$ws.Cells.Item(31,5).Value = '  +7.73%  '

# Row 32
$ws.Cells.Item(32,4).Value = '''6.352'
$ws.Cells.Item(32,5).Value = '  +20.89%  '

# Row 33
$ws.Cells.Item(33,4).Value = '''4.035'
$ws.Cells.Item(33,5).Value = '  +0.14%  '

# Row 34
$ws.Cells.Item(34,4).Value = '''0.9959'
$ws.Cells.Item(34,5).Value = '  +16.48%  '

# Row 35
$ws.Cells.Item(35,4).Value = '''1.758'
$ws.Cells.Item(35,5).Value = '  +16.62%  '

# Row 36
$ws.Cells.Item(36,4).Value = '''0.08449'
$ws.Cells.Item(36,5).Value = '  +5.29%  '

# Row 37
$ws.Cells.Item(37,4).Value = '''12.46'
$ws.Cells.Item(37,5).Value = '  +15.79%  '

# Row 38
$ws.Cells.Item(38,2).Value = 'Hedera'
$ws.Cells.Item(38,3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(38,4).Value = '''0.06450'
$ws.Cells.Item(38,5).Value = '  +9.51%  '

# Row 39
$ws.Cells.Item(39,2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(39,3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(39,4).Value = '''5.370'
$ws.Cells.Item(39,5).Value = '  +7.58%  '

# Row 40
$ws.Cells.Item(40,2).Value = 'FraxShare'
$ws.Cells.Item(40,3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(40,4).Value = '''8.904'
$ws.Cells.Item(40,5).Value = '  +14.86%  '

# Row 41
$ws.Cells.Item(41,4).Value = '''1.297'
$ws.Cells.Item(41,5).Value = '  +6.79%  '

# Row 42
$ws.Cells.Item(42,4).Value = '''0.02350'
$ws.Cells.Item(42,5).Value = '  +11.51%  '

# Row 43
$ws.Cells.Item(43,4).Value = '''0.2119'
$ws.Cells.Item(43,5).Value = '  +10.31%  '

# Row 44
$ws.Cells.Item(44,4).Value = '''0.6155'
$ws.Cells.Item(44,5).Value = '  +12.38%  '

# Row 45
$ws.Cells.Item(45,4).Value = '''0.9948'
$ws.Cells.Item(45,5).Value = '  +3.70%  '

# Row 46
$ws.Cells.Item(46,2).Value = 'PancakeSwap'
$ws.Cells.Item(46,3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(46,4).Value = '''3.798'
$ws.Cells.Item(46,5).Value = '  +5.97%  '

# Row 47
$ws.Cells.Item(47,2).Value = 'EnergySwap'
$ws.Cells.Item(47,3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(47,4).Value = '''13.19'
$ws.Cells.Item(47,5).Value = '  +5.62%  '

# Row 48
$ws.Cells.Item(48,4).Value = '''0.5999'
$ws.Cells.Item(48,5).Value = '  +9.62%  '

# Row 49
$ws.Cells.Item(49,4).Value = '''127.84'
$ws.Cells.Item(49,5).Value = '  +4.81%  '

# Row 50
$ws.Cells.Item(50,4).Value = '''2.025'
$ws.Cells.Item(50,5).Value = '  +7.83%  '

# Row 51
$ws.Cells.Item(51,4).Value = '''0.07149'
$ws.Cells.Item(51,5).Value = '  +7.71%  '
